# Generate Report for Handback
# Updates the "a983145c-02b6-4462-a5e2-85fcb89a3b19" row (row 7) on both the
# zh-cn and de-de localization-status sheets now that a handback has been
# received: fill in the Latest Target File / Latest Handback File / Latest
# Handback DateTime columns (and, for de-de, an Error Detail message about
# the handback not being based on the very latest source revision).

$wb = $excel.ActiveWorkbook

$targetFileName = "a983145c-02b6-4462-a5e2-85fcb89a3b19.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37c59eb000304debc6fc6bd51546fff76775a805/e2e/a983145c-02b6-4462-a5e2-85fcb89a3b19.md"

function Set-HandbackRow($sheetName, $handbackFile, $handbackDateTime, $errorDetail) {

    $ws = $wb.Worksheets.Item($sheetName)

    # I7 - Latest Target File: becomes a hyperlink to the source .md file,
    # mirroring the look/behaviour of the A7 cell.
    $i7 = $ws.Range("I7")
    $i7.Value = $targetFileName
    $f = $i7.Font()
    $f.Underline = 2
    $f.Color = 15570276
    $ws.Hyperlinks.Add($i7, $targetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetFileName)

    # J7 - Latest Handback File
    $ws.Range("J7").Value = $handbackFile

    # K7 - Latest Handback DateTime
    $ws.Range("K7").Value = $handbackDateTime

    if ($errorDetail) {
        # P7 - Error Detail
        $ws.Range("P7").Value = $errorDetail
    }
}

Set-HandbackRow `
    "zh-cn" `
    "a983145c-02b6-4462-a5e2-85fcb89a3b19.83bd9c1c5e0b47501dcac27568e32f6bb929b552.zh-cn.xlf" `
    "2016-09-03 05:00:36" `
    ""

Set-HandbackRow `
    "de-de" `
    "a983145c-02b6-4462-a5e2-85fcb89a3b19.83bd9c1c5e0b47501dcac27568e32f6bb929b552.de-de.xlf" `
    "2016-09-03 05:00:43" `
    "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ffaf5c838e26db46d060600224bc1ee8e4b0fa48/e2e/a983145c-02b6-4462-a5e2-85fcb89a3b19.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37c59eb000304debc6fc6bd51546fff76775a805/e2e/a983145c-02b6-4462-a5e2-85fcb89a3b19.md."

Write-Host "Handback report updated."
